$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header in E1
$ws.Range("E1").Value = "forebrain_neurons"

# Fill E2:E65 with the constant value for all data rows
for ($row = 2; $row -le 65; $row++) {
    $ws.Cells.Item($row, 5).Value = 24560000000
}
